# Regenerate the Bmp7-Acvr2b LR-pair TPM table with the updated NATMI output.
# A new sending cluster ("ECs") is added and all TPM-derived metrics for the
# Bmp7 -> Acvr2b edges are recomputed, growing the data block from 8 to 12 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@("ECs", "Bmp7", "Acvr2b", "ECs", 1, 0.3333333333333333, 0.03443933333333333, 0.103318, 0.05823261822459219, 0.0582326182245922, 3, 1, 1.716657, 5.149971, 0.3840886036988016, 0.3840886036988015, 0.05912052264199999, 0.532084703778, 0.022366485023609, 0.022366485023609)
    ,@("ECs", "Bmp7", "Acvr2b", "FAPs", 1, 0.3333333333333333, 0.03443933333333333, 0.103318, 0.05823261822459219, 0.0582326182245922, 3, 1, 0.3517506666666667, 1.055252, 0.07870146593648156, 0.07870146593648154, 0.01211405845955556, 0.109026526136, 0.004582992419594877, 0.004582992419594877)
    ,@("ECs", "Bmp7", "Acvr2b", "MuSCs", 1, 0.3333333333333333, 0.03443933333333333, 0.103318, 0.05823261822459219, 0.0582326182245922, 3, 1, 1.677572333333333, 5.032717, 0.3753437146230962, 0.3753437146230962, 0.05777447277844443, 0.5199702550059999, 0.02185724723664705, 0.02185724723664705)
    ,@("ECs", "Bmp7", "Acvr2b", "Resolving-Mac", 1, 0.3333333333333333, 0.03443933333333333, 0.103318, 0.05823261822459219, 0.0582326182245922, 3, 1, 0.7234496666666667, 2.170349, 0.1618662157416207, 0.1618662157416207, 0.02491512422022222, 0.224236117982, 0.009425893544741276, 0.009425893544741276)
    ,@("FAPs", "Bmp7", "Acvr2b", "ECs", 2, 0.6666666666666666, 0.51625, 1.54875, 0.8729143757654733, 0.8729143757654734, 3, 1, 1.716657, 5.149971, 0.3840886036988016, 0.3840886036988015, 0.8862241762499999, 7.97601758625, 0.3352764637363717, 0.3352764637363717)
    ,@("FAPs", "Bmp7", "Acvr2b", "FAPs", 2, 0.6666666666666666, 0.51625, 1.54875, 0.8729143757654733, 0.8729143757654734, 3, 1, 0.3517506666666667, 1.055252, 0.07870146593648156, 0.07870146593648154, 0.1815912816666667, 1.634321535, 0.06869964100977145, 0.06869964100977145)
    ,@("FAPs", "Bmp7", "Acvr2b", "MuSCs", 2, 0.6666666666666666, 0.51625, 1.54875, 0.8729143757654733, 0.8729143757654734, 3, 1, 1.677572333333333, 5.032717, 0.3753437146230962, 0.3753437146230962, 0.8660467170833333, 7.79442045375, 0.327642924347714, 0.327642924347714)
    ,@("FAPs", "Bmp7", "Acvr2b", "Resolving-Mac", 2, 0.6666666666666666, 0.51625, 1.54875, 0.8729143757654733, 0.8729143757654734, 3, 1, 0.7234496666666667, 2.170349, 0.1618662157416207, 0.1618662157416207, 0.3734808904166667, 3.36132801375, 0.1412953466716163, 0.1412953466716163)
    ,@("MuSCs", "Bmp7", "Acvr2b", "ECs", 1, 0.3333333333333333, 0.04072033333333334, 0.122161, 0.06885300600993445, 0.06885300600993445, 3, 1, 1.716657, 5.149971, 0.3840886036988016, 0.3840886036988015, 0.06990284525900001, 0.629125607331, 0.02644565493882092, 0.02644565493882092)
    ,@("MuSCs", "Bmp7", "Acvr2b", "FAPs", 1, 0.3333333333333333, 0.04072033333333334, 0.122161, 0.06885300600993445, 0.06885300600993445, 3, 1, 0.3517506666666667, 1.055252, 0.07870146593648156, 0.07870146593648154, 0.01432340439688889, 0.128910639572, 0.005418832507115216, 0.005418832507115215)
    ,@("MuSCs", "Bmp7", "Acvr2b", "MuSCs", 1, 0.3333333333333333, 0.04072033333333334, 0.122161, 0.06885300600993445, 0.06885300600993445, 3, 1, 1.677572333333333, 5.032717, 0.3753437146230962, 0.3753437146230962, 0.06831130460411111, 0.6148017414370001, 0.02584354303873517, 0.02584354303873516)
    ,@("MuSCs", "Bmp7", "Acvr2b", "Resolving-Mac", 1, 0.3333333333333333, 0.04072033333333334, 0.122161, 0.06885300600993445, 0.06885300600993445, 3, 1, 0.7234496666666667, 2.170349, 0.1618662157416207, 0.1618662157416207, 0.02945911157655556, 0.265132004189, 0.01114497552526316, 0.01114497552526316)
)

$numRows = $rows.Count
$numCols = 20
$data = New-Object "object[,]" $numRows, $numCols
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $data[$i, $j] = $rows[$i][$j]
    }
}

# Data block now spans A2:T13 (previously A2:T9).
$ws.Range("A2:T13").Value = $data
